{"js": "// Fix the typo \"controlle\" -> \"controle\" (\"zegel controlle\" -> \"zegel controle\")\n// inside the legacy FORMTEXT form-field answer text near the top of the document.\nconst body = context.document.body;\nconst results = body.search(\"controlle\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"controle\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Fix the typo \"controlle\" -> \"controle\" (\"zegel controlle\" -> \"zegel controle\")\n# inside the legacy FORMTEXT form-field answer text near the top of the document.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"controlle\"\n$find.Replacement.Text = \"controle\"\n$find.Execute(\n    \"controlle\",  # FindText\n    $false,       # MatchCase\n    $false,       # MatchWholeWord\n    $false,       # MatchWildcards\n    $false,       # MatchSoundsLike\n    $false,       # MatchAllWordForms\n    $true,        # Forward\n    1,            # Wrap (wdFindContinue)\n    $false,       # Format\n    \"controle\",   # ReplaceWith\n    2             # Replace (wdReplaceAll)\n)\n"}
